$d = $word.ActiveDocument

# Remove leading space before "Norte"
$d.Content.Find.Execute(" Norte", $true, $false, $false, $false, $false,
                         $true, 1, $false, "Norte", 2)

# Remove leading space before "Este"
$d.Content.Find.Execute(" Este", $true, $false, $false, $false, $false,
                         $true, 1, $false, "Este", 2)

# Add a space after the dollar sign for $3,800.00
$d.Content.Find.Execute("$3,800.00", $true, $false, $false, $false, $false,
                         $true, 1, $false, "$ 3,800.00", 2)

# Add a space after the dollar sign for $1,200.00
$d.Content.Find.Execute("$1,200.00", $true, $false, $false, $false, $false,
                         $true, 1, $false, "$ 1,200.00", 2)
